# Update 9 juni - list utang
# Applies the edits described by the commit: clears the "Tanggal Pembayaran"
# dates and one stray mandor name from rows 126-132, fixes two misspelled
# "Perorangan" entries, normalizes the date styling on rows 127-132 back to
# the sheet's usual style, and appends two new payment rows (133-134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the payment-date (column A) cells for rows 126-132 ---
# These rows no longer have a "Tanggal Pembayaran" value.
$ws.Range("A126:A132").Clear()

# --- 2. Clear the stray mandor name in B128 ---
$ws.Range("B128").Clear()

# --- 3. Fix the misspelled "Perorangan" entries ---
# B127 was "perorngan", B131 was "perorangan" -> both become "Perorangan"
$ws.Range("B127").Value = "Perorangan"
$ws.Range("B131").Value = "Perorangan"

# --- 4. Normalize the number formatting on D127:E132 back to the style ---
# --- used by the rest of the table (copy format only from row 126)    ---
$ws.Range("D126").Copy()
$ws.Range("D127:D132").PasteSpecial(-4122)
$ws.Range("E126").Copy()
$ws.Range("E127:E132").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Append the two new payment rows ---
$ws.Range("A133").Value = 45450
$ws.Range("B133").Value = "Aceng"
$ws.Range("D133").Value = 45423
$ws.Range("E133").Value = 45436
$ws.Range("F133").Value = 5049000
$ws.Range("G133").Value = 3448000
$ws.Range("I133").Formula = "=G133-F133"
$ws.Range("J133").Value = 3448000
$ws.Range("K133").Formula = "=G133-J133"
$ws.Range("L133").Formula = "=G133-J133+H133"
$ws.Range("M133").Value = 100000

$ws.Range("B134").Value = "Haji Skun"
$ws.Range("D134").Value = 45422
$ws.Range("E134").Value = 45435
$ws.Range("F134").Value = 16778000
$ws.Range("G134").Value = 16778000
$ws.Range("H134").Value = 3523000
$ws.Range("I134").Formula = "=G134-F134"
$ws.Range("J134").Value = 18000000
$ws.Range("K134").Formula = "=G134-J134"
$ws.Range("L134").Formula = "=G134-J134+H134"
$ws.Range("M134").Value = 100000

# --- 6. Expand Table1 / AutoFilter to include the two new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:N134"))

# --- 7. Update the selection / scroll position like the saved workbook ---
$ws.Range("H137").Select()
